$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")

# Row 20
$ws.Range("H20").Value = 952.3333
$ws.Range("I20").Value = 446.375
$ws.Range("J20").Value = 5000
$ws.Range("K20").Value = 446.375
$ws.Range("L20").Value = 5000
$ws.Range("M20").Value = -216.375
$ws.Range("N20").Value = -5460

# Row 34
$ws.Range("H34").Value = 5604
$ws.Range("I34").Value = 5604
$ws.Range("K34").Value = 5604
$ws.Range("M34").Value = -5401

# Row 35
$ws.Range("H35").Value = 952.3333
$ws.Range("I35").Value = 446.375
$ws.Range("J35").Value = 5000
$ws.Range("K35").Value = 446.375
$ws.Range("L35").Value = 5000
$ws.Range("M35").Value = -67.375
$ws.Range("N35").Value = -5758

# Row 36
$ws.Range("H36").Value = 5604
$ws.Range("I36").Value = 5604
$ws.Range("K36").Value = 5604
$ws.Range("M36").Value = -4889

# Row 46
$ws.Range("H46").Value = 4998.5
$ws.Range("I46").Value = 4998.5
$ws.Range("K46").Value = 14995.5
$ws.Range("M46").Value = -14876.5

# Row 49
$ws.Range("H49").Value = 975
$ws.Range("J49").Value = 975
$ws.Range("L49").Value = 2925
$ws.Range("N49").Value = -3197

# Row 59
$ws.Range("H59").Value = 5000
$ws.Range("J59").Value = 5000
$ws.Range("L59").Value = 15000
$ws.Range("N59").Value = -16114

# Row 60
$ws.Range("H60").Value = 4998.5
$ws.Range("I60").Value = 4998.5
$ws.Range("K60").Value = 14995.5
$ws.Range("M60").Value = -14511.5

# Row 76
$ws.Range("H76").Value = 4999.25
$ws.Range("I76").Value = 4999
$ws.Range("K76").Value = 4999
$ws.Range("M76").Value = -4684

# Row 79
$ws.Range("H79").Value = 4999.25
$ws.Range("I79").Value = 4999
$ws.Range("K79").Value = 4999
$ws.Range("M79").Value = -3907

# Row 92
$ws.Range("H92").Value = 740.7143
$ws.Range("I92").Value = 740.7143
$ws.Range("K92").Value = 740.7143
$ws.Range("M92").Value = 507.2857

# Row 96
$ws.Range("H96").Value = 529.8
$ws.Range("I96").Value = 412.5
$ws.Range("J96").Value = 999
$ws.Range("K96").Value = 1237.5
$ws.Range("L96").Value = 2997
$ws.Range("M96").Value = 135.5
$ws.Range("N96").Value = -5743


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")

# Row 2
$ws.Range("H2").Value = 3458.5715
$ws.Range("I2").Value = 837.4
$ws.Range("J2").Value = 10011.5
$ws.Range("K2").Value = 837.4
$ws.Range("L2").Value = 10011.5
$ws.Range("M2").Value = -724.4
$ws.Range("N2").Value = -10237.5

# Row 39
$ws.Range("H39").Value = 2435.5715
$ws.Range("I39").Value = 2006.2
$ws.Range("J39").Value = 3509
$ws.Range("K39").Value = 2006.2
$ws.Range("L39").Value = 3509
$ws.Range("M39").Value = -1486.2
$ws.Range("N39").Value = -4549

# Row 50
$ws.Range("H50").Value = 28316.334
$ws.Range("I50").Value = 20499.5
$ws.Range("J50").Value = 43950
$ws.Range("K50").Value = 20499.5
$ws.Range("L50").Value = 43950
$ws.Range("M50").Value = -19785.5
$ws.Range("N50").Value = -45378

# Row 56
$ws.Range("H56").Value = 25000
$ws.Range("I56").Value = 25000
$ws.Range("K56").Value = 25000
$ws.Range("M56").Value = -24258

# Row 102
$ws.Range("H102").Value = 2384.9285
$ws.Range("I102").Value = 932.2222
$ws.Range("J102").Value = 4999.8
$ws.Range("K102").Value = 932.2222
$ws.Range("L102").Value = 4999.8
$ws.Range("M102").Value = 689.7778
$ws.Range("N102").Value = -8243.799999999999

# Row 116
$ws.Range("H116").Value = 3458.5715
$ws.Range("I116").Value = 837.4
$ws.Range("J116").Value = 10011.5
$ws.Range("K116").Value = 837.4
$ws.Range("L116").Value = 10011.5
$ws.Range("M116").Value = 1456.6
$ws.Range("N116").Value = -14599.5


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")

# Row 3
$ws.Range("H3").Value = 3458.5715
$ws.Range("I3").Value = 837.4
$ws.Range("J3").Value = 10011.5
$ws.Range("K3").Value = 837.4
$ws.Range("L3").Value = 10011.5
$ws.Range("M3").Value = -723.4
$ws.Range("N3").Value = -10239.5

# Row 20
$ws.Range("H20").Value = 1933.875

# Row 32
$ws.Range("H32").Value = 25000
$ws.Range("J32").Value = 25000
$ws.Range("L32").Value = 25000
$ws.Range("N32").Value = -25768

# Row 36
$ws.Range("H36").Value = 2184.25
$ws.Range("I36").Value = 2184.25
$ws.Range("K36").Value = 2184.25
$ws.Range("M36").Value = -1650.25

# Row 38
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

# Row 99
$ws.Range("H99").Value = 2538.0527
$ws.Range("I99").Value = 2412.4285
$ws.Range("K99").Value = 2412.4285
$ws.Range("M99").Value = -914.4285


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")

# Row 16
$ws.Range("H16").Value = 951.5714
$ws.Range("J16").Value = 983.3333
$ws.Range("L16").Value = 983.3333
$ws.Range("N16").Value = -1557.3333

# Row 33
$ws.Range("H33").Value = 3715.75
$ws.Range("I33").Value = 3715.75
$ws.Range("K33").Value = 3715.75
$ws.Range("M33").Value = -3336.75

# Row 105
$ws.Range("H105").Value = 1753.4546
$ws.Range("I105").Value = 1037.4
$ws.Range("J105").Value = 2350.1667
$ws.Range("K105").Value = 1037.4
$ws.Range("L105").Value = 2350.1667
$ws.Range("M105").Value = 709.5999999999999
$ws.Range("N105").Value = -5844.1667

# Row 113
$ws.Range("H113").Value = 951.5714
$ws.Range("J113").Value = 983.3333
$ws.Range("L113").Value = 983.3333
$ws.Range("N113").Value = -5323.3333


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")

# Row 33
$ws.Range("H33").Value = 2998
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 2998
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 17988
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -18554

# Row 57
$ws.Range("H57").Value = 2550
$ws.Range("I57").Value = 2550
$ws.Range("K57").Value = 7650
$ws.Range("M57").Value = -7091

# Row 68
$ws.Range("H68").Value = 1037.5
$ws.Range("I68").Value = 1037.5
$ws.Range("K68").Value = 3112.5
$ws.Range("M68").Value = -2301.5

# Row 71
$ws.Range("H71").Value = 1037.5
$ws.Range("I71").Value = 1037.5
$ws.Range("K71").Value = 9337.5
$ws.Range("M71").Value = -5281.5


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")

# Row 113
$ws.Range("H113").Value = 859.5714
$ws.Range("I113").Value = 859.5714
$ws.Range("K113").Value = 859.5714
$ws.Range("M113").Value = 1310.4286


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")

# Row 4
$ws.Range("H4").Value = 1169.6666
$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 1750
$ws.Range("K4").Value = 9
$ws.Range("L4").Value = 1750
$ws.Range("M4").Value = 104
$ws.Range("N4").Value = -1976

# Row 28
$ws.Range("H28").Value = 1169.6666
$ws.Range("I28").Value = 9
$ws.Range("J28").Value = 1750
$ws.Range("K28").Value = 9
$ws.Range("L28").Value = 1750
$ws.Range("M28").Value = 223
$ws.Range("N28").Value = -2214

# Row 29
$ws.Range("H29").Value = 20016
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

# Row 31
$ws.Range("H31").Value = 4248.5713
$ws.Range("I31").Value = 4248.5713
$ws.Range("K31").Value = 4248.5713
$ws.Range("M31").Value = -4000.5713

# Row 34
$ws.Range("H34").Value = 14666.333
$ws.Range("J34").Value = 14499.5
$ws.Range("L34").Value = 14499.5
$ws.Range("N34").Value = -14843.5

# Row 35
$ws.Range("H35").Value = 836.375
$ws.Range("I35").Value = 836.375
$ws.Range("K35").Value = 836.375
$ws.Range("M35").Value = -500.375

# Row 37
$ws.Range("H37").Value = 1169.6666
$ws.Range("I37").Value = 9
$ws.Range("J37").Value = 1750
$ws.Range("K37").Value = 9
$ws.Range("L37").Value = 1750
$ws.Range("M37").Value = 98
$ws.Range("N37").Value = -1964

# Row 58
$ws.Range("H58").Value = 9725.5
$ws.Range("I58").Value = 8599.666999999999
$ws.Range("J58").Value = 13103
$ws.Range("K58").Value = 8599.666999999999
$ws.Range("L58").Value = 13103
$ws.Range("M58").Value = -8339.666999999999
$ws.Range("N58").Value = -13623

# Row 61
$ws.Range("H61").Value = 2454
$ws.Range("I61").Value = 2268.8572
$ws.Range("K61").Value = 2268.8572
$ws.Range("M61").Value = -2066.8572

# Row 93
$ws.Range("H93").Value = 998
$ws.Range("J93").Value = 799.5
$ws.Range("L93").Value = 799.5
$ws.Range("N93").Value = -3295.5

# Row 113
$ws.Range("H113").Value = 2454
$ws.Range("I113").Value = 2268.8572
$ws.Range("K113").Value = 2268.8572
$ws.Range("M113").Value = -98.85719999999992

# Row 138
$ws.Range("H138").Value = 80000
$ws.Range("J138").Value = 80000
$ws.Range("L138").Value = 80000
$ws.Range("N138").Value = -90280


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")

# Row 32
$ws.Range("H32").Value = 14750
$ws.Range("J32").Value = 14500
$ws.Range("L32").Value = 14500
$ws.Range("N32").Value = -15134

# Row 55
$ws.Range("H55").Value = 13450.333
$ws.Range("I55").Value = 149
$ws.Range("K55").Value = 149
$ws.Range("M55").Value = 128

# Row 113
$ws.Range("H113").Value = 898
$ws.Range("I113").Value = 898
$ws.Range("K113").Value = 2694
$ws.Range("M113").Value = -524

# Row 136
$ws.Range("H136").Value = 2998
$ws.Range("I136").Value = 2797.9
$ws.Range("K136").Value = 8393.700000000001
$ws.Range("M136").Value = -5843.700000000001

# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# Row 141
$ws.Range("H141").Value = 89999
$ws.Range("I141").Value = 89999
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 89999
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -84819
$ws.Range("N141").ClearContents()

